$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: continue the running "index" shared-formula pattern
# (=A{row-1}+1) that already fills A2:A420, down through row 551.
# Each new batch is written as a single leading cell followed by the
# remaining multi-cell range, which keeps the generated shared-formula
# groups well formed (each with a proper master/ref cell).
$ws.Range("A421").Formula = "=A420+1"
$ws.Range("A422:A453").Formula = "=A421+1"

$ws.Range("A454").Formula = "=A453+1"
$ws.Range("A455:A517").Formula = "=A454+1"

$ws.Range("A518").Formula = "=A517+1"
$ws.Range("A519:A551").Formula = "=A518+1"

# Column C: occurrence count of 1 for every newly-indexed row (421-551).
$ws.Range("C421:C551").Value = 1

# Update the saved selection/active cell to match the new end of the list,
# and scroll the view down so row 546 is back at the top (mirrors the
# author scrolling down after appending the new rows).
$ws.Range("B557").Select()
$excel.ActiveWindow.ScrollRow = 546
$excel.ActiveWindow.ScrollColumn = 1
